# Update classifier predictions to evaluate accuracy
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("D37").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0
